# "Generate Report for Archive"
# - Flip every "Ready for handoff" status cell to "In Translation"
#   (Overview!E2/F2 - the zh-cn/de-de status roll-up columns - and the
#   corresponding Status column on each per-locale sheet).
# - Narrow the "Status" column (it no longer needs to fit the long
#   "Ready for handoff" label) on the Overview sheet (cols E & F) and on
#   each locale sheet (col C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status values everywhere they appear.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Shrink the now-narrower Status columns.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
